$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "fantasy points" column (G),
# shifting it to column I, to make room for "height" and "weight".
$ws.Range("G1:H1").EntireColumn.Insert()

# New header cells
$ws.Range("G1").Value = "height"
$ws.Range("H1").Value = "weight"

# Copy the style of an existing header cell (e.g. F1) onto the new headers
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in constant height/weight values for each data row (rows 2-17)
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 7).Value = 6.333333333333333   # G: height
    $ws.Cells.Item($r, 8).Value = 238                  # H: weight
}
